# Tidsplan.xlsx - "Rettede tidsplan for dagen."
#
# The day's coding-task cells (F12:F15) get reshuffled:
#   F12 "Kodning af "lavere klasser""                  -> "Diskussion af controller-\nimplementation"
#   F13 "Diskussion af controller-\nimplementation"     -> "Diskussion af GUI-\nmockups"
#   F14 "Implementation af \ntoString()-metoder..."     -> "Kodning af \n"lavere klasser""  (reworded)
#   F15 "Diskussion af GUI-\nmockups"                   -> "Reservation Hashset/MapSet \nmellem pers. og sæde" (new)
#
# F12 additionally picks up the word-wrap formatting already used by F13-F15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing texts up one slot …
$ws.Range("F12").Value = "Diskussion af controller-`nimplementation"
$ws.Range("F13").Value = "Diskussion af GUI-`nmockups"

# … and fill in the brand-new / reworded items (F15 first so the shared-string
# table gets "Reservation …" at index 42 and "Kodning …" at index 43, matching
# the order the strings were typed in on the day).
$ws.Range("F15").Value = "Reservation Hashset/MapSet `nmellem pers. og sæde"
$ws.Range("F14").Value = "Kodning af `n`"lavere klasser`""

# F12 used to have a non-wrapping style (s="49"); match the wrapped style
# used by the rest of the column (s="65") now that it holds a two-line text.
$ws.Range("F12").WrapText = $true

# Restore the cursor position recorded in the saved file.
$ws.Range("G15").Select() | Out-Null
